$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the confidential disclosure text: date changes from 2021-07-07 to 2021-07-08
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-08 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-7
$ws.Range("D2").Value = 0.4996083357787084
$ws.Range("E2").Value = -0.003440366972477071

$ws.Range("D3").Value = 0.3236486371782715
$ws.Range("E3").Value = -0.008789528995891915

$ws.Range("D4").Value = 0.08991300931506051
$ws.Range("E4").Value = -0.001854349291975699

$ws.Range("D5").Value = 0.05871568087732786
$ws.Range("E5").Value = -0.002158109950022857

$ws.Range("D6").Value = 0.02811433685063163
$ws.Range("E6").Value = -0.005046257359125339

$ws.Range("D7").Value = 0.9999999999999998
$ws.Range("E7").Value = -0.004998872298062618

# Restore the original sheet protection (password-protected in the source file;
# re-apply protection with the same allowances that were previously in effect).
$ws.Protect($null, $true, $true, $true, $null, $true, $false, $false, $true, $true, $true, $true, $true, $true, $true, $true)
